# Update Rainfall 4 sheet with water_plus_bottle_mass_collected(g) values
# pulled from volume data (B3 R4), rows 38-55, column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    38 = 827
    39 = 827
    40 = 893
    41 = 856
    42 = 817
    43 = 817
    44 = 896
    45 = 808
    46 = 808
    47 = 846
    48 = 903
    49 = 868
    50 = 915
    51 = 863
    52 = 878
    53 = 883
    54 = 899
    55 = 875
}

foreach ($row in $values.Keys) {
    $ws.Range("H$row").Value = $values[$row]
}

# Move the selection to match the saved view state after the edit.
[void]$ws.Range("Q52").Select()
